$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = "2025-05-05 11:18:12"

# Row 3
$ws.Range("B3").Value = "NA"
# C3 holds the text "0" (was "4"); force text so Excel doesn't coerce it
# to a number, then drop the formatting we had to apply so the cell's
# style stays the same as before (no explicit style, like its neighbors).
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "0"
$ws.Range("C3").ClearFormats()
$ws.Range("D3").Value = "NA"
$ws.Range("E3").Value = "NA"
$ws.Range("F3").Value = "NA"
$ws.Range("G3").Value = "NA"
$ws.Range("H3").Value = "2025-05-05 11:19:41"
$ws.Range("I3").Value = "Error de consulta o no registrado"

# Row 4
$ws.Range("B4").Value = "NA"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "0"
$ws.Range("C4").ClearFormats()
$ws.Range("D4").Value = "NA"
$ws.Range("E4").Value = "NA"
$ws.Range("F4").Value = "NA"
$ws.Range("G4").Value = "NA"
$ws.Range("H4").Value = "2025-05-05 11:21:11"
$ws.Range("I4").Value = "Error de consulta o no registrado"
